$wb = $excel.ActiveWorkbook

# The sheet holding "Gesamtinvestitionskosten" (total investment cost) data
$ws = $wb.Worksheets.Item("Gesamtinvestitionskosten")
$ws.Activate()

# Update the input quantities in column B (rows 2-10)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 7
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 9

# Update the percentage parameters used by the formulas above
$ws.Range("B20").Value = 0.12
$ws.Range("B21").Value = 0.11

# Update the selected range shown when the sheet is reopened
$ws.Range("A1:F2").Select()

$wb.Save()
